# Apply the crypto price/volume refresh for Sat Sep  9 23:32:17 UTC 2023 run.
# D = Price column, E = Volume(1h) column (row N corresponds to worksheet row N).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '25.938.31'
$ws.Cells.Item(2, 5).Value = '  -0.06%  '

$ws.Cells.Item(3, 4).Value = '1.637.74'
$ws.Cells.Item(3, 5).Value = '  -0.16%  '

$ws.Cells.Item(5, 4).Value = '''214.56'
$ws.Cells.Item(5, 5).Value = '  -0.19%  '

$ws.Cells.Item(6, 5).Value = '  -0.20%  '

$ws.Cells.Item(7, 5).Value = '  -0.03%  '

$ws.Cells.Item(8, 5).Value = '  -0.45%  '

$ws.Cells.Item(9, 4).Value = '''0.0636'
$ws.Cells.Item(9, 5).Value = '  -0.42%  '

$ws.Cells.Item(10, 4).Value = '''19.51'
$ws.Cells.Item(10, 5).Value = '  -0.81%  '

$ws.Cells.Item(11, 5).Value = '  +0.06%  '

$ws.Cells.Item(12, 5).Value = '  -0.27%  '

$ws.Cells.Item(13, 4).Value = '1.604.24'
$ws.Cells.Item(13, 5).Value = '  -0.62%  '

$ws.Cells.Item(14, 4).Value = '''0.542'
$ws.Cells.Item(14, 5).Value = '  -0.51%  '

$ws.Cells.Item(15, 4).Value = '''63.22'
$ws.Cells.Item(15, 5).Value = '  +0.90%  '

$ws.Cells.Item(16, 4).Value = '0.0₃0758'
$ws.Cells.Item(16, 5).Value = '  -0.57%  '

$ws.Cells.Item(17, 4).Value = '25.960.31'
$ws.Cells.Item(17, 5).Value = '  -0.02%  '

$ws.Cells.Item(18, 5).Value = '  -0.03%  '

$ws.Cells.Item(19, 4).Value = '''193.90'
$ws.Cells.Item(19, 5).Value = '  -0.21%  '

$ws.Cells.Item(20, 5).Value = '  -1.03%  '

$ws.Cells.Item(21, 5).Value = '  -0.81%  '

$ws.Cells.Item(22, 4).Value = '''6.18'
$ws.Cells.Item(22, 5).Value = '  -1.69%  '

$ws.Cells.Item(23, 5).Value = '  +3.73%  '

$ws.Cells.Item(24, 4).Value = '''143.61'
$ws.Cells.Item(24, 5).Value = '  -0.36%  '

$ws.Cells.Item(25, 4).Value = '''1.00'
$ws.Cells.Item(25, 5).Value = '  -0.15%  '

$ws.Cells.Item(26, 5).Value = '  -0.67%  '

$ws.Cells.Item(27, 5).Value = '  +0.24%  '

$ws.Cells.Item(28, 4).Value = '''15.47'
$ws.Cells.Item(28, 5).Value = '  -0.12%  '

$ws.Cells.Item(29, 5).Value = '  -0.05%  '

$ws.Cells.Item(30, 4).Value = '''0.0494'
$ws.Cells.Item(30, 5).Value = '  -1.39%  '

$ws.Cells.Item(31, 5).Value = '  -0.96%  '

$ws.Cells.Item(33, 5).Value = '  -0.78%  '

$ws.Cells.Item(34, 5).Value = '  +0.62%  '

$ws.Cells.Item(35, 4).Value = '''0.900'
$ws.Cells.Item(35, 5).Value = '  -0.52%  '

$ws.Cells.Item(36, 4).Value = '1.127.30'
$ws.Cells.Item(36, 5).Value = '  -1.17%  '

$ws.Cells.Item(37, 5).Value = '  -1.49%  '

$ws.Cells.Item(38, 5).Value = '  -0.38%  '

$ws.Cells.Item(39, 5).Value = '  -0.75%  '

$ws.Cells.Item(40, 4).Value = '''98.51'
$ws.Cells.Item(40, 5).Value = '  -0.98%  '

$ws.Cells.Item(41, 5).Value = '  -0.17%  '

$ws.Cells.Item(42, 4).Value = '''0.792'
$ws.Cells.Item(42, 5).Value = '  -1.07%  '

$ws.Cells.Item(43, 5).Value = '  -0.39%  '

$ws.Cells.Item(44, 4).Value = '''56.30'
$ws.Cells.Item(44, 5).Value = '  -0.56%  '

$ws.Cells.Item(45, 5).Value = '  +2.20%  '

$ws.Cells.Item(46, 5).Value = '  -1.61%  '

$ws.Cells.Item(47, 4).Value = '''7.71'
$ws.Cells.Item(47, 5).Value = '  +1.07%  '

$ws.Cells.Item(48, 5).Value = '  -0.54%  '

$ws.Cells.Item(49, 5).Value = '  -0.03%  '

$ws.Cells.Item(50, 5).Value = '  -2.13%  '

$ws.Cells.Item(51, 4).Value = '''5.50'
$ws.Cells.Item(51, 5).Value = '  -0.60%  '
